$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the "_old" / "_new" suffixed header labels (row 1) to the
#    version-qualified "_FV2210" / "_FV2304" labels. Column K ("diff") is
#    left untouched.
# ---------------------------------------------------------------------------
$fv2210Headers = @("Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210","Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210")
$fv2304Headers = @("Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304","Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the used range A1:U81 into an Excel Table ("Table1") with an
#    auto filter on the header row.
# ---------------------------------------------------------------------------
$range = $ws.Range("A1:U81")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.TableStyle = ""
